# DG: update model API
# Locate the sequence-diagram label shape that reads "deletePerson(p)" and
# update the model API call to "deletePersons(p)".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame) {
        if ($candidate.TextFrame.TextRange.Text -eq "deletePerson(p)") {
            $targetShape = $candidate
            break
        }
    }
}

if ($targetShape -eq $null) {
    # Fallback: known shape name from the original deck.
    $targetShape = $s.Shapes.Item("TextBox 77")
}

$tr = $targetShape.TextFrame.TextRange

# "deletePerson" -> "deletePersons"
$nameRun = $tr.Characters(1, 12)
$nameRun.Text = "deletePersons"

# "(p)" -> split into "(p" and ")" runs (matches target OOXML run structure)
$closeParen = $tr.Characters($tr.Length, 1)
$closeParen.Font.Color.RGB = $closeParen.Font.Color.RGB
